$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:E values per row (F column / "Win" is unchanged).
$data = @{
    2  = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987)
    3  = @(0.1169995834814548, 0.3048912486333797, 18.71679738969934, 0.5333859586016987)
    4  = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 13.86384647080068)
    5  = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987)
    6  = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987)
    7  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987)
    8  = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 13.86384647080068)
    9  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987)
    10 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987)
    11 = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987)
    12 = @(1.445647641019636, 9.983522426115931, 3.223369029078222, 13.86384647080068)
    13 = @(1.445647641019636, 9.983522426115931, 3.223369029078222, 13.86384647080068)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]

    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 7).Value = ($b + $c + $d + $e)
}
